$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.288.22"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.490.34"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "2.880.42"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "2.491.02"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "47.187.64"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +14.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.10%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "1.995.57"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.20%  "
